# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, reflecting freshly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 228
$ws1.Range("F6").Value = 1601
$ws1.Range("F8").Value = 423
$ws1.Range("F10").Value = 212
$ws1.Range("F11").Value = 6806
$ws1.Range("F13").Value = 537
$ws1.Range("F14").Value = 116
$ws1.Range("F16").Value = 2372
$ws1.Range("F20").Value = 97
$ws1.Range("F27").Value = 117

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 70
$ws2.Range("F3").Value = 19

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 679

# --- Sheet 4: 全部类型 (All types, aggregate of the above) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 679
$ws4.Range("F7").Value = 70
$ws4.Range("F8").Value = 228
$ws4.Range("F9").Value = 1601
$ws4.Range("F10").Value = 19
$ws4.Range("F13").Value = 423
$ws4.Range("F15").Value = 212
$ws4.Range("F16").Value = 6806
$ws4.Range("F18").Value = 537
$ws4.Range("F19").Value = 116
$ws4.Range("F21").Value = 2372
$ws4.Range("F25").Value = 97
$ws4.Range("F32").Value = 117
